$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(36, 1).Value = 2943549
$ws.Cells.Item(36, 2).Value = "Pril Power blue 650 ml"
$ws.Cells.Item(36, 3).Value = 12
$ws.Cells.Item(36, 4).Value = 164

$ws.Range("A36:D36").Select()
